{"js": "// Word JS API (Office.js) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1) \"Ativa\u00e7\u00e3o: 01/01/2012\" -> \"Ativa\u00e7\u00e3o: 01/01/2023\"\n// ---------------------------------------------------------------------\nconst dateHits = body.search(\"Ativa\u00e7\u00e3o: 01/01/2012\", { matchCase: true });\ndateHits.load(\"items\");\nawait context.sync();\nif (dateHits.items.length > 0) {\n  dateHits.items[0].insertText(\"Ativa\u00e7\u00e3o: 01/01/2023\", \"Replace\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2) After the \"Objetivos\" body paragraph, insert the italic English\n//    translation as a brand-new paragraph.\n// ---------------------------------------------------------------------\nconst objetivosHits = body.search(\n  \"Fornecer ao aluno o conhecimento das principais t\u00e9cnicas de caracteriza\u00e7\u00e3o f\u00edsica e qu\u00edmica de materiais.\",\n  { matchCase: true }\n);\nobjetivosHits.load(\"items\");\nawait context.sync();\nif (objetivosHits.items.length > 0) {\n  const objetivosPara = objetivosHits.items[0].paragraphs.getFirst();\n  const objetivosEn = objetivosPara.insertParagraph(\n    \"Provide the student with knowledge of the main techniques of physical and chemical characterization of materials.\",\n    \"After\"\n  );\n  objetivosEn.font.italic = true;\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3) After the \"Programa resumido\" body paragraph, insert the italic\n//    English translation as a brand-new paragraph.\n// ---------------------------------------------------------------------\nconst resumidoHits = body.search(\n  \"An\u00e1lise granulom\u00e9trica e superficial. An\u00e1lises microestruturais. An\u00e1lises t\u00e9rmicas. Reometria.\",\n  { matchCase: true }\n);\nresumidoHits.load(\"items\");\nawait context.sync();\nif (resumidoHits.items.length > 0) {\n  const resumidoPara = resumidoHits.items[0].paragraphs.getFirst();\n  const resumidoEn = resumidoPara.insertParagraph(\n    \"Granulometric and surface analysis. Microstructural analyses. Thermal analysis. Rheometry.\",\n    \"After\"\n  );\n  resumidoEn.font.italic = true;\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 4) The long \"Programa\" paragraph: collapse its four runs (joined by\n//    manual line breaks) into a single run with no breaks.\n// ---------------------------------------------------------------------\nconst programaMerged =\n  \"An\u00e1lise granulom\u00e9trica. Adsor\u00e7\u00e3o BET, porosidade e picnometria. \" +\n  \"An\u00e1lises microestruturais: difra\u00e7\u00e3o de raios X, figura de Laue; espalhamento de raios X (SAXS). Difra\u00e7\u00e3o de el\u00e9trons. Microscopia \u00d3ptica. Microscopia eletr\u00f4nica, microan\u00e1lise de raios X (EDX e WDX). \" +\n  \"An\u00e1lises t\u00e9rmicas: An\u00e1lise t\u00e9rmica diferencial (DTA), calorimetria explorat\u00f3ria diferencial (DSC) e termogravimetria (TGA).\" +\n  \"Reometria de l\u00edquidos, solu\u00e7\u00f5es e pastas.\";\n\nconst programaHits = body.search(\"Reometria de l\u00edquidos, solu\u00e7\u00f5es e pastas.\", {\n  matchCase: true\n});\nprogramaHits.load(\"items\");\nawait context.sync();\nif (programaHits.items.length > 0) {\n  const programaPara = programaHits.items[0].paragraphs.getFirst();\n  const programaRange = programaPara.getRange();\n  programaRange.insertText(programaMerged, \"Replace\");\n  await context.sync();\n\n  // -------------------------------------------------------------------\n  // 5) Right after the (now single-run) \"Programa\" paragraph, insert the\n  //    italic English translation as a brand-new paragraph.\n  // -------------------------------------------------------------------\n  const programaEn = programaPara.insertParagraph(\n    \"Grain size analysis. BET adsorption, porosity and pycnometry.\" +\n      \"Microstructural analysis: X-ray diffraction, Laue figure; X-ray scattering (SAXS). Electron diffraction. Optical Microscopy. Electron microscopy, X-ray microanalysis (EDX and WDX).\" +\n      \"Thermal analysis: Differential thermal analysis (DTA), differential scanning calorimetry (DSC) and thermogravimetry (TGA).\" +\n      \"Rheometry of liquids, solutions and pastes.\",\n    \"After\"\n  );\n  programaEn.font.italic = true;\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument / $d resolve to the open document.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($doc, $text) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        if ($doc.Paragraphs.Item($i).Range.Text.TrimEnd() -eq $text) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# ---------------------------------------------------------------------\n# 1) \"Ativa\u00e7\u00e3o: 01/01/2012\" -> \"Ativa\u00e7\u00e3o: 01/01/2023\"\n# ---------------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n[void]$find.Execute(\"Ativa\u00e7\u00e3o: 01/01/2012\", $false, $false, $false, $false, $false, $true, 1, $false, \"Ativa\u00e7\u00e3o: 01/01/2023\", 2)\n\n# ---------------------------------------------------------------------\n# 2) After the \"Objetivos\" body paragraph, insert the italic English\n#    translation as a brand-new paragraph.\n# ---------------------------------------------------------------------\n$idx = Find-ParagraphIndex $d \"Fornecer ao aluno o conhecimento das principais t\u00e9cnicas de caracteriza\u00e7\u00e3o f\u00edsica e qu\u00edmica de materiais.\"\nif ($idx -gt 0) {\n    $p = $d.Paragraphs.Item($idx)\n    $p.Range.InsertParagraphAfter()\n    $newp = $d.Paragraphs.Item($idx + 1)\n    $newp.Range.Text = \"Provide the student with knowledge of the main techniques of physical and chemical characterization of materials.\"\n    $r = $newp.Range\n    [void]$r.MoveEnd(1, -1)\n    $r.Font.Italic = $true\n}\n\n# ---------------------------------------------------------------------\n# 3) After the \"Programa resumido\" body paragraph, insert the italic\n#    English translation as a brand-new paragraph.\n# ---------------------------------------------------------------------\n$idx = Find-ParagraphIndex $d \"An\u00e1lise granulom\u00e9trica e superficial. An\u00e1lises microestruturais. An\u00e1lises t\u00e9rmicas. Reometria.\"\nif ($idx -gt 0) {\n    $p = $d.Paragraphs.Item($idx)\n    $p.Range.InsertParagraphAfter()\n    $newp = $d.Paragraphs.Item($idx + 1)\n    $newp.Range.Text = \"Granulometric and surface analysis. Microstructural analyses. Thermal analysis. Rheometry.\"\n    $r = $newp.Range\n    [void]$r.MoveEnd(1, -1)\n    $r.Font.Italic = $true\n}\n\n# ---------------------------------------------------------------------\n# 4) The long \"Programa\" paragraph: collapse its four runs (joined by\n#    manual line breaks) into a single run with no breaks.\n# ---------------------------------------------------------------------\n$programaMerged = \"An\u00e1lise granulom\u00e9trica. Adsor\u00e7\u00e3o BET, porosidade e picnometria. \" + `\n    \"An\u00e1lises microestruturais: difra\u00e7\u00e3o de raios X, figura de Laue; espalhamento de raios X (SAXS). Difra\u00e7\u00e3o de el\u00e9trons. Microscopia \u00d3ptica. Microscopia eletr\u00f4nica, microan\u00e1lise de raios X (EDX e WDX). \" + `\n    \"An\u00e1lises t\u00e9rmicas: An\u00e1lise t\u00e9rmica diferencial (DTA), calorimetria explorat\u00f3ria diferencial (DSC) e termogravimetria (TGA).\" + `\n    \"Reometria de l\u00edquidos, solu\u00e7\u00f5es e pastas.\"\n\n$idx = Find-ParagraphIndex $d \"Reometria de l\u00edquidos, solu\u00e7\u00f5es e pastas.\"\nif ($idx -lt 0) {\n    # Fallback: locate by the paragraph's distinctive opening text in case\n    # the trailing-text probe above doesn't match (defensive only).\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        if ($d.Paragraphs.Item($i).Range.Text.StartsWith(\"An\u00e1lise granulom\u00e9trica. Adsor\u00e7\u00e3o BET\")) {\n            $idx = $i\n            break\n        }\n    }\n}\nif ($idx -gt 0) {\n    $p = $d.Paragraphs.Item($idx)\n    $p.Range.Text = $programaMerged\n\n    # ---------------------------------------------------------------\n    # 5) Right after the (now single-run) \"Programa\" paragraph, insert\n    #    the italic English translation as a brand-new paragraph.\n    # ---------------------------------------------------------------\n    $p.Range.InsertParagraphAfter()\n    $newp = $d.Paragraphs.Item($idx + 1)\n    $programaEn = \"Grain size analysis. BET adsorption, porosity and pycnometry.\" + `\n        \"Microstructural analysis: X-ray diffraction, Laue figure; X-ray scattering (SAXS). Electron diffraction. Optical Microscopy. Electron microscopy, X-ray microanalysis (EDX and WDX).\" + `\n        \"Thermal analysis: Differential thermal analysis (DTA), differential scanning calorimetry (DSC) and thermogravimetry (TGA).\" + `\n        \"Rheometry of liquids, solutions and pastes.\"\n    $newp.Range.Text = $programaEn\n    $r = $newp.Range\n    [void]$r.MoveEnd(1, -1)\n    $r.Font.Italic = $true\n}\n"}
